# Trade #100 closed at 2026-02-16 21:40:03 - leadlag DOWN +0.000%
#
# This script updates:
#  - Summary sheet aggregate stats (rows 2 & 3)
#  - leadlag sheet: trades #64 & #65 (rows 53 & 54) transition OPEN -> CLOSED,
#    and a brand new trade #100 (row 76) is appended as OPEN
#  - All Trades sheet: appends the now-CLOSED trades #64 & #65 as rows 65 & 66
#  - Comparison sheet aggregate stats for the leadlag strategy (row 2)

$wb = $excel.ActiveWorkbook

function Set-TextCell {
    param($ws, [string]$addr, [string]$val)
    # Force text storage so Excel doesn't auto-coerce dates/times/percents/numeric-looking
    # strings into numbers or serial dates.
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $val
}

function Set-NumCell {
    param($ws, [string]$addr, $val)
    $ws.Range($addr).Value = $val
}

# ---------------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")

Set-NumCell  $wsSummary "C2" 65
Set-TextCell $wsSummary "D2" "66.2%"
Set-TextCell $wsSummary "E2" "+15.9421%"
Set-TextCell $wsSummary "F2" "+0.2453%"

Set-NumCell  $wsSummary "C3" 74
Set-TextCell $wsSummary "D3" "44.6%"
Set-TextCell $wsSummary "E3" "+11.0913%"
Set-TextCell $wsSummary "F3" "+0.1499%"

# ---------------------------------------------------------------------------
# leadlag sheet
# ---------------------------------------------------------------------------
$wsLeadlag = $wb.Worksheets.Item("leadlag")

# Trade #64 (row 53): OPEN -> CLOSED
Set-NumCell  $wsLeadlag "G53" 68778.37411600001
Set-TextCell $wsLeadlag "H53" "CLOSED"
Set-NumCell  $wsLeadlag "I53" -0.1302
Set-NumCell  $wsLeadlag "J53" -1.3
Set-TextCell $wsLeadlag "M53" "time_exit_5min"
Set-NumCell  $wsLeadlag "N53" 5

# Trade #65 (row 54): OPEN -> CLOSED
Set-NumCell  $wsLeadlag "G54" 68420.180271
Set-TextCell $wsLeadlag "H54" "CLOSED"
Set-NumCell  $wsLeadlag "I54" 0.2314
Set-NumCell  $wsLeadlag "J54" 2.31
Set-TextCell $wsLeadlag "M54" "time_exit_5min"
Set-NumCell  $wsLeadlag "N54" 5

# New trade #100 (row 76): freshly opened
Set-NumCell  $wsLeadlag "A76" 100
Set-TextCell $wsLeadlag "B76" "2026-02-16"
Set-TextCell $wsLeadlag "C76" "21:40:03"
Set-TextCell $wsLeadlag "D76" "leadlag"
Set-TextCell $wsLeadlag "E76" "DOWN"
Set-NumCell  $wsLeadlag "F76" 68276.55
Set-TextCell $wsLeadlag "H76" "OPEN"
Set-NumCell  $wsLeadlag "I76" 0
Set-NumCell  $wsLeadlag "J76" 0
Set-NumCell  $wsLeadlag "K76" 0.75
Set-TextCell $wsLeadlag "L76" "Binance leading with -0.081% move"
Set-NumCell  $wsLeadlag "N76" 0

# ---------------------------------------------------------------------------
# All Trades sheet: append the two newly-closed trades
# ---------------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("All Trades")

# Row 65 <- leadlag trade #64 (now closed)
Set-NumCell  $wsAll "A65" 64
Set-TextCell $wsAll "B65" "2026-02-16"
Set-TextCell $wsAll "C65" "21:34:53"
Set-TextCell $wsAll "D65" "leadlag"
Set-TextCell $wsAll "E65" "DOWN"
Set-NumCell  $wsAll "F65" 68688.94
Set-NumCell  $wsAll "G65" 68778.37411600001
Set-TextCell $wsAll "H65" "CLOSED"
Set-NumCell  $wsAll "I65" -0.1302
Set-NumCell  $wsAll "J65" -1.3
Set-NumCell  $wsAll "K65" 0.75
Set-TextCell $wsAll "L65" "Binance leading with -0.082% move"
Set-TextCell $wsAll "M65" "time_exit_5min"
Set-NumCell  $wsAll "N65" 5

# Row 66 <- leadlag trade #65 (now closed)
Set-NumCell  $wsAll "A66" 65
Set-TextCell $wsAll "B66" "2026-02-16"
Set-TextCell $wsAll "C66" "21:34:59"
Set-TextCell $wsAll "D66" "leadlag"
Set-TextCell $wsAll "E66" "DOWN"
Set-NumCell  $wsAll "F66" 68578.89999999999
Set-NumCell  $wsAll "G66" 68420.180271
Set-TextCell $wsAll "H66" "CLOSED"
Set-NumCell  $wsAll "I66" 0.2314
Set-NumCell  $wsAll "J66" 2.31
Set-NumCell  $wsAll "K66" 0.75
Set-TextCell $wsAll "L66" "Binance leading with -0.216% move"
Set-TextCell $wsAll "M66" "time_exit_5min"
Set-NumCell  $wsAll "N66" 5

# ---------------------------------------------------------------------------
# Comparison sheet
# ---------------------------------------------------------------------------
$wsComparison = $wb.Worksheets.Item("Comparison")

Set-NumCell  $wsComparison "B2" 74
Set-TextCell $wsComparison "C2" "44.6%"
Set-TextCell $wsComparison "D2" "2.83"
Set-TextCell $wsComparison "E2" "+0.5196%"
Set-TextCell $wsComparison "F2" "-0.3028%"
Set-TextCell $wsComparison "G2" "1.72"
